# Added additional CCS scenarios
# Edits the "connections" sheet:
#   - C14: simple_BF -> all
#   - Appends 5 new rows (25-29) describing biofuel / fossil-fuel inflow scenarios

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("connections")

# Existing row tweak
$ws.Cells.Item(14, 3).Value = "all"

# New row 25: lime | simple_coke | inflow | biofuel | fuel | outflows | simple_fuel | fuel
$ws.Cells.Item(25, 2).Value = "lime"
$ws.Cells.Item(25, 3).Value = "simple_coke"
$ws.Cells.Item(25, 4).Value = "inflow"
$ws.Cells.Item(25, 5).Value = "biofuel"
$ws.Cells.Item(25, 6).Value = "fuel"
$ws.Cells.Item(25, 6).Font.Color = 0
$ws.Cells.Item(25, 7).Value = "outflows"
$ws.Cells.Item(25, 8).Value = "simple_fuel"
$ws.Cells.Item(25, 8).NumberFormat = "@"
$ws.Cells.Item(25, 9).Value = "fuel"
$ws.Cells.Item(25, 9).Font.Color = 0

# New row 26: steel | simple_sinter | inflow | fossil fuel | fuel | outflows | simple_fuel | fuel
$ws.Cells.Item(26, 2).Value = "steel"
$ws.Cells.Item(26, 3).Value = "simple_sinter"
$ws.Cells.Item(26, 4).Value = "inflow"
$ws.Cells.Item(26, 5).Value = "fossil fuel"
$ws.Cells.Item(26, 6).Value = "fuel"
$ws.Cells.Item(26, 7).Value = "outflows"
$ws.Cells.Item(26, 8).Value = "simple_fuel"
$ws.Cells.Item(26, 8).NumberFormat = "@"
$ws.Cells.Item(26, 9).Value = "fuel"
$ws.Cells.Item(26, 9).Font.Color = 0

# New row 27: steel | simple_sinter | inflow | biofuel | fuel | outflows | simple_fuel | fuel
$ws.Cells.Item(27, 2).Value = "steel"
$ws.Cells.Item(27, 3).Value = "simple_sinter"
$ws.Cells.Item(27, 4).Value = "inflow"
$ws.Cells.Item(27, 5).Value = "biofuel"
$ws.Cells.Item(27, 6).Value = "fuel"
$ws.Cells.Item(27, 7).Value = "outflows"
$ws.Cells.Item(27, 8).Value = "simple_fuel"
$ws.Cells.Item(27, 8).NumberFormat = "@"
$ws.Cells.Item(27, 9).Value = "fuel"
$ws.Cells.Item(27, 9).Font.Color = 0

# New row 28: pellets | simple_pellets | inflow | biofuel | fuel | outflows | simple_fuel | fuel
$ws.Cells.Item(28, 2).Value = "pellets"
$ws.Cells.Item(28, 3).Value = "simple_pellets"
$ws.Cells.Item(28, 4).Value = "inflow"
$ws.Cells.Item(28, 5).Value = "biofuel"
$ws.Cells.Item(28, 6).Value = "fuel"
$ws.Cells.Item(28, 7).Value = "outflows"
$ws.Cells.Item(28, 8).Value = "simple_fuel"
$ws.Cells.Item(28, 8).NumberFormat = "@"
$ws.Cells.Item(28, 9).Value = "fuel"
$ws.Cells.Item(28, 9).Font.Color = 0

# New row 29: pellets | simple_pellets | inflow | fossil fuel | fuel | outflows | simple_fuel | fuel
$ws.Cells.Item(29, 2).Value = "pellets"
$ws.Cells.Item(29, 3).Value = "simple_pellets"
$ws.Cells.Item(29, 4).Value = "inflow"
$ws.Cells.Item(29, 5).Value = "fossil fuel"
$ws.Cells.Item(29, 6).Value = "fuel"
$ws.Cells.Item(29, 7).Value = "outflows"
$ws.Cells.Item(29, 8).Value = "simple_fuel"
$ws.Cells.Item(29, 8).NumberFormat = "@"
$ws.Cells.Item(29, 9).Value = "fuel"
$ws.Cells.Item(29, 9).Font.Color = 0

# Update selection to match the author's final view (cosmetic, matches diff)
$ws.Range("C33").Select()
